$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.459.25"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.30%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.604.21"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.23%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.06"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.517"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +7.04%  "
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "26.86"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +10.64%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "43.49"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.67%  "
$ws.Range("E10").Value = "  +2.75%  "
$ws.Range("E11").Value = "  +2.68%  "
$ws.Range("E12").Value = "  +1.99%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.832.06"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.15%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.606.15"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.36%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "29.496.59"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.53%  "
$ws.Range("E16").Value = "  +4.96%  "
$ws.Range("E17").Value = "  +3.00%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.48"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.76%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "240.36"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.94%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.62"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.53%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0692"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.20%  "
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("E23").Value = "  +3.28%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.15"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.88%  "
$ws.Range("E25").Value = "  -0.35%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.41"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.27"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.63%  "
$ws.Range("E28").Value = "  +5.18%  "
$ws.Range("E29").Value = "  +2.21%  "
$ws.Range("E31").Value = "  +2.81%  "
$ws.Range("E32").Value = "  +0.11%  "
$ws.Range("E33").Value = "  +2.28%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.422.88"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.38%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.09"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.37%  "
$ws.Range("E36").Value = "  -1.27%  "
$ws.Range("E37").Value = "  +1.96%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.83"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.04%  "
$ws.Range("E39").Value = "  +0.03%  "
$ws.Range("E40").Value = "  +2.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.533"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.09%  "
$ws.Range("E42").Value = "  +0.81%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "53.07"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +21.98%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.793"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.19%  "
$ws.Range("E46").Value = "  +1.77%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "65.26"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.54%  "
$ws.Range("E48").Value = "  -0.22%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.744.25"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "86.42"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.46%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.835"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.95%  "
